$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is incremented by 1 day
# (45560 -> 45561) for every data row (rows 2 through 29).
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45560) {
        $cell.Value2 = 45561
    }
}
